$d = $word.ActiveDocument

# 1. Title: "ตาราง … Use case Description" -> "ตารางที่ 1 Use case Description"
$d.Content.Find.Execute("ตาราง …", $true, $false, $false, $false, $false, $true, 1, $false, "ตารางที่ 1", 2)

# 2. Table cell: merge "1.3.1" text
$d.Content.Find.Execute("1.3.1", $true, $false, $false, $false, $false, $true, 1, $false, "1.3.1", 2)
